# ---------------------------------------------------------------------------
# "description for pipeline behaviour"
#
# Renames Sheet2 -> "pipeline behaviour" and fills it in with a description
# of the functional-unit / pipeline-stage timeline for `lw r6, 0(r6)`,
# mirroring the layout/styling already used on the RAW sheet. Also restores
# the RAW sheet's selection to "select all" and makes the new sheet active.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$raw = $wb.Worksheets.Item("RAW")
$ws  = $wb.Worksheets.Item("Sheet2")

$ws.Name = "pipeline behaviour"

# --- column widths (mirror RAW's col A, widen B:F for the longer labels) ---
$ws.Columns.Item(1).ColumnWidth = $raw.Columns.Item(1).ColumnWidth
$ws.Range("B:F").ColumnWidth = 13.8
$ws.Range("G:AR").ColumnWidth = $raw.Columns.Item(7).ColumnWidth

# --- A1: title ---
$raw.Range("A1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$ws.Range("A1").Value = "Assuming a 5-stage pipeline"

# --- A3: "Pipeline Stage" ---
$raw.Range("A3").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Value = "Pipeline Stage"

# --- row 4: instruction + pipeline-stage header row ---
$raw.Range("A3:F3").Copy() | Out-Null
$ws.Range("A4:F4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = "lw r6, 0(r6)"
$ws.Range("B4").Value = "IF"
$ws.Range("C4").Value = "ID"
$ws.Range("D4").Value = "EX"
$ws.Range("E4").Value = "MEM"
$ws.Range("F4").Value = "WB"

# --- A6: "Functional Unit" ---
$raw.Range("A3").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").Value = "Functional Unit"

# --- row 7: functional-unit header cells (IM / Reg / ALU / IM / Reg) ---
$raw.Range("C4:D4").Copy() | Out-Null
$ws.Range("B7:C7").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7:F7").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Value = "IM"
$ws.Range("C7").Value = "Reg"
$ws.Range("D7").Value = "ALU"
$ws.Range("E7").Value = "IM"
$ws.Range("F7").Value = "Reg"

# --- row 8: behaviour description cells ---
$raw.Range("B4").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").Value = "lw"

$raw.Range("C4:F4").Copy() | Out-Null
$ws.Range("C8:F8").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").Value = "reg[r6]"
$ws.Range("D8").Value = "addr = 0 + reg[r6]"
$ws.Range("E8").Value = "mem[addr]"
$ws.Range("F8").Value = "reg[r6] = mem[addr]"

$wb.Application.CutCopyMode = $false

# --- selection / active sheet bookkeeping ---
$raw.Cells.Select() | Out-Null
$ws.Select() | Out-Null
$ws.Range("C21").Select() | Out-Null
